$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1569408081499378
$ws.Range("H2").Value = 19.41859865002044
$ws.Range("G3").Value = 0.1717004433804159
$ws.Range("H3").Value = 92.84125121145496
$ws.Range("G4").Value = -0.6422546796980082
$ws.Range("H4").Value = -4.373369765375687
$ws.Range("G5").Value = -0.5589303733638036
$ws.Range("H5").Value = 8.461094492249089
$ws.Range("G6").Value = 0.2121750077080789
$ws.Range("H6").Value = -13.7877702720971
$ws.Range("G7").Value = 0.3841434140024276
$ws.Range("H7").Value = 134.4883765612052
$ws.Range("G8").Value = 0.1118222900081235
$ws.Range("H8").Value = -32.35556536435304
$ws.Range("G9").Value = 0.2135692896278156
$ws.Range("H9").Value = 9.470762354544124
$ws.Range("G10").Value = -0.1948932169126427
$ws.Range("H10").Value = -241.0422975155985
$ws.Range("G11").Value = -0.1431703458337011
$ws.Range("H11").Value = -20.54590162976265
$ws.Range("G12").Value = 0.1775940895285582
$ws.Range("H12").Value = 11.67382701988623
$ws.Range("G13").Value = 0.1818157329640047
$ws.Range("H13").Value = -11.59477377251484
$ws.Range("G14").Value = 0.1557858288607722
$ws.Range("H14").Value = -17.7360204018228
$ws.Range("G15").Value = 0.2485130745104423
$ws.Range("H15").Value = -0.5571480287567659
$ws.Range("G16").Value = -0.007929441062452365
$ws.Range("H16").Value = -121.736259311429
$ws.Range("G17").Value = 0.01336789697259665
$ws.Range("H17").Value = -62.31273019544194
$ws.Range("G18").Value = 0.1243236406798293
$ws.Range("H18").Value = -28.2669141303163
$ws.Range("G19").Value = 0.1688672931566585
$ws.Range("H19").Value = 34.28496222913127
$ws.Range("G20").Value = 0.1195450499911382
$ws.Range("H20").Value = 4.26773990028963
$ws.Range("G21").Value = 0.1141721250941084
$ws.Range("H21").Value = 13.71858587358322
$ws.Range("G22").Value = 0.04705371038395056
$ws.Range("H22").Value = -50.04761736526709
$ws.Range("G23").Value = 0.07499123292960645
$ws.Range("H23").Value = -30.87727779881895
$ws.Range("G24").Value = -0.2514461386044626
$ws.Range("H24").Value = -101.7632540314836
$ws.Range("G25").Value = -0.1572611063694041
$ws.Range("H25").Value = 29.3062632897336
$ws.Range("G26").Value = 0.1693744719191226
$ws.Range("H26").Value = 6.537008041640522
$ws.Range("G27").Value = 0.1959373634708947
$ws.Range("H27").Value = -2.259218691758295
$ws.Range("G28").Value = 0.02272725869068091
$ws.Range("H28").Value = 383.0135060249247
$ws.Range("G29").Value = 0.0004393102680317901
$ws.Range("H29").Value = -97.143160455098
$ws.Range("I2").Value = -3.175617967186845
